# Update the Virginia school-reopening tracker:
#   - Shorten/lowercase the header row labels
#   - Correct the "Cumberland County" date of return (was mis-keyed to
#     1/7/2020, should be 1/7/2021)
#   - Leave the selection on the corrected cell, scrolled back to the top
#     of the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: County/September/November/Date of return/Link/Notes (not
# comprehensive) -> shorter lowercase labels.
$ws.Range("A1").Value = "county"
$ws.Range("B1").Value = "september"
$ws.Range("C1").Value = "november"
$ws.Range("D1").Value = "date"
$ws.Range("E1").Value = "link"
$ws.Range("F1").Value = "notes"

# Fix the Cumberland County return date (row 8): 1/7/2020 -> 1/7/2021.
$ws.Range("D8").Value = (Get-Date -Year 2021 -Month 1 -Day 7 -Hour 0 -Minute 0 -Second 0)

# Bring the view back to the top and leave the corrected cell selected.
[void]$ws.Range("D8").Select()

"done"
